$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PSSM values for B2:K21 (rows 2-21, columns B-K), one row per array element
$newValues = @(
    @(-18.67493704156459, 1.993315191769324, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459),
    @(-18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, 0.6844543566876048, -18.67493704156459, -18.67493704156459),
    @(-18.67493704156459, 1.965711258700977, 1.647784431876333, -18.67493704156459, 3.419641013307861, -18.67493704156459, 4.321924820105167, -18.67493704156459, 0.7332123012826389, -18.67493704156459),
    @(-18.67493704156459, 1.663278123764789, -18.67493704156459, -18.67493704156459, -18.67493704156459, 3.075793248670317, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459),
    @(-18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459),
    @(2.448169543790034, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459),
    @(-18.67493704156459, -18.67493704156459, -18.67493704156459, 1.85108548691833, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459),
    @(3.862225656815108, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459),
    @(-18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, 1.54252018576372, -18.67493704156459, 2.112446656377669),
    @(-18.67493704156459, -18.67493704156459, -18.67493704156459, 2.919521062777622, -18.67493704156459, 2.634522103608808, -18.67493704156459, -18.67493704156459, -18.67493704156459, 1.852670728770675),
    @(-18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459),
    @(-18.67493704156459, -18.67493704156459, -18.67493704156459, 2.485247688128342, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, 1.715707412638179, 1.509678046252768),
    @(-18.67493704156459, -18.67493704156459, 1.522069624053639, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, 2.149139907526764),
    @(-18.67493704156459, -18.67493704156459, 1.753432332806205, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459),
    @(-18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, 1.914491921524432, -18.67493704156459),
    @(-18.67493704156459, 2.15651943392733, 1.852632550645412, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, 2.401944432255072, 2.507993726932036, -18.67493704156459),
    @(-18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, 2.168914253821159, 2.484126794733818, -18.67493704156459),
    @(-18.67493704156459, -18.67493704156459, 2.046203044200311, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459, 1.857568230790926, -18.67493704156459, -18.67493704156459),
    @(-18.67493704156459, 1.02630341721077, 1.526599429272496, -18.67493704156459, 3.217106539876137, -18.67493704156459, -18.67493704156459, 1.052333652534309, -18.67493704156459, 2.257295735903721),
    @(-18.67493704156459, 1.292707176619843, -18.67493704156459, 1.69019033044248, -18.67493704156459, 2.421894550060511, -18.67493704156459, -18.67493704156459, -18.67493704156459, -18.67493704156459)
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $rowValues = $newValues[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $colNum = $j + 2
        $ws.Cells.Item($rowNum, $colNum).Value = $rowValues[$j]
    }
}

